$d = $word.ActiveDocument

# Insert a new paragraph right after the existing last paragraph.
# InsertParagraphAfter() mirrors the formatting (pPr/rPr, incl. lang)
# of the paragraph it is called on, matching the target markup.
$lastPara = $d.Paragraphs.Last.Range
$lastPara.InsertParagraphAfter()

# Populate the newly created (now last) paragraph with its text.
$d.Paragraphs.Last.Range.Text = "2 – Random order event"
